# Daily attendance processing - 2026-01-27 01:47:03
# For every row in the "Recorded By" column whose value lists the
# recorders as "System, <email>", flip the order to "<email>, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Locate the "Recorded By" column dynamically from the header row (row 1).
$recordedByCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item(1, $c).Value()
    if ($header -eq "Recorded By") {
        $recordedByCol = $c
        break
    }
}

if ($recordedByCol -eq 0) {
    $recordedByCol = 7
}

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"
$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Value()
    if ($val -eq $oldValue) {
        $cell.Value = $newValue
        $changed = $changed + 1
    }
}

Write-Output "Updated Recorded By order on $changed row(s)."
